$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = "{'cluster': 15, 'training': 1330, 'window': 200, 'threshold': 8.0}"
$ws.Range("J2").Value = 0.5077524951193482

# Row 3
$ws.Range("G3").Value = "[1600, 1700, 1800, 4000, 4100, 4200, 4400, 4500, 4700, 4800, 4900]"
$ws.Range("I3").Value = "{'cluster': 14, 'training': 1256, 'window': 200, 'threshold': 2.5}"
$ws.Range("J3").Value = 0.7068190251011401

# Row 4
$ws.Range("I4").Value = "{'cluster': 16, 'training': 1236, 'window': 200, 'threshold': 4.5}"
$ws.Range("J4").Value = 0.8544361060485244

# Row 5
$ws.Range("I5").Value = "{'cluster': 17, 'training': 751, 'window': 200, 'threshold': 2.5}"
$ws.Range("J5").Value = 0.5201737459283322

# Row 6
$ws.Range("I6").Value = "{'cluster': 19, 'training': 440, 'window': 200, 'threshold': 3.0}"
$ws.Range("J6").Value = 0.4550081070046872

# Row 7
$ws.Range("G7").Value = "[10500, 10800]"
$ws.Range("H7").Value = 0.4762962962962962
$ws.Range("I7").Value = "{'cluster': 10, 'training': 4479, 'window': 600, 'threshold': 3.5}"
$ws.Range("J7").Value = 0.4726196799892932

# Row 8
$ws.Range("G8").Value = "[340, 360, 380, 740, 820, 1080, 1140, 1240, 1260, 1540, 1660, 1860, 1960, 2040, 2560, 2580, 2880, 3160, 3440, 4000, 4060, 4140, 4160, 4360, 4500, 4800, 5280, 5340, 5420, 5640, 5660, 5760, 5940, 6240, 6440, 6520, 6640, 6740, 6840, 6860, 6880, 6900, 6920, 7040, 7460, 7520, 7640, 7920, 7960, 8560, 8780, 8860, 8940, 9200, 9300, 9320, 9340, 9440, 9740, 9760, 10040, 10060, 10180, 10200, 10460, 10480, 11040, 11060, 11400, 11660, 11680, 11880, 11940, 12300, 12480, 12660, 12720, 12920, 12960, 12980, 13100, 13120, 13160, 13560, 13620, 13760, 13880, 14100, 14120, 14140, 14220, 14260, 14280, 14360, 14400, 14420, 14480, 14780, 14900]"
$ws.Range("H8").Value = 0.01704545454545455
$ws.Range("I8").Value = "{'cluster': 10, 'training': 222, 'window': 40, 'threshold': 4.5}"
$ws.Range("J8").Value = 4.840415807906538

# Row 9
$ws.Range("G9").Value = "[2300, 4800]"
$ws.Range("H9").Value = 0.5
$ws.Range("I9").Value = "{'cluster': 15, 'training': 2915, 'window': 200, 'threshold': 1.5}"
$ws.Range("J9").Value = 2.022873396985233

# Row 10
$ws.Range("G10").Value = "[300, 375, 2100, 2175, 2250, 2325, 2400]"
$ws.Range("H10").Value = 0.3636363636363636
$ws.Range("I10").Value = "{'cluster': 13, 'training': 178, 'window': 150, 'threshold': 1.0}"
$ws.Range("J10").Value = 2.639962041983381

# Row 11
$ws.Range("G11").Value = "[5000]"
$ws.Range("H11").Value = 0.6666666666666666
$ws.Range("I11").Value = "{'cluster': 14, 'training': 1011, 'window': 100, 'threshold': 4.5}"
$ws.Range("J11").Value = 3.780867037829012

# Row 13
$ws.Range("G13").Value = "[4392, 4758, 4941, 6588, 6771, 12810]"
$ws.Range("H13").Value = 0.2857142857142858
$ws.Range("I13").Value = "{'cluster': 18, 'training': 4210, 'window': 366, 'threshold': 8.5}"
$ws.Range("J13").Value = 0.8104959439951926
